$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = '[Janice L%Atkins%NULL%0, Jane A H%Masoli%NULL%2, Jane A H%Masoli%NULL%0, Joao%Delgado%NULL%2, Joao%Delgado%NULL%0, Luke C%Pilling%NULL%2, Luke C%Pilling%NULL%0, Chia-Ling%Kuo%NULL%2, Chia-Ling%Kuo%NULL%0, George A%Kuchel%NULL%2, George A%Kuchel%NULL%0, David%Melzer%D.Melzer@exeter.ac.uk%1, Anne B%Newman%NULL%2, Anne B%Newman%NULL%0]'
$ws.Range("I2").Value2 = ''
$ws.Range("J2").Value2 = 'Oxford University Press'
$ws.Range("E3").Value2 = '[Angelo%Bianchetti%angelo.bianchetti@grupposandonato.it%0, R.%Rozzini%NULL%1, F.%Guerini%NULL%1, S.%Boffelli%NULL%1, P.%Ranieri%NULL%1, G.%Minelli%NULL%1, L.%Bianchetti%NULL%1, M.%Trabucchi%NULL%1]'
$ws.Range("I3").Value2 = ''
$ws.Range("J3").Value2 = 'Springer Paris'
$ws.Range("D4").Value2 = 'Our study aimed to investigate the factors affecting the prognosis of patients with mental disorders with COVID-19. All patients with mental disorders who were diagnosed with COVID-19 at the intensive care unit of Wuhan Mental Health Center during the period January 3 to March 1, 2020 were selected.
 The influence of the baseline characteristics, clinical symptoms, laboratory parameters and the types of mental disorders on prognosis were analyzed.
 According to their final prognosis, the patients were divided into the deceased group (5 patients) and the cured group (25 patients).
 The mortality rate of patients with dementia was significantly higher than that of patients with other mental disorders (P = 0.001).
 The levels of certain laboratory parameters in the serum of dementia patients were significantly increased compared with levels in nondementia patients (WBC count: 10.100±6.147 vs.
 5.694±3.383, p = 0.029; neutrophil count: 8.504± 5.993 vs.
 3.764 ± 2.733, P = 0.008; BUN: 8.300± 4.072 vs.
 4.364 ± 1.196, P = 0.001).
 Our research indicated that the mortality rate of dementia patients with COVID-19 was higher than that of patients with other mental disorders.
 A focus on the inflammatory response of dementia patients may provide novel ideas for reducing mortality.
'
$ws.Range("E4").Value2 = '[Yan%Wan%NULL%0, Juan%Wu%NULL%2, Lihua%Ni%NULL%1, Qinqin%Luo%NULL%1, Cheng%Yuan%NULL%1, Fang%Fan%NULL%1, Hong%Liu%NULL%0, Changjiang%Zhang%NULL%2, Yuandi%Xiang%NULL%1, Qin%Xie%NULL%1]'
$ws.Range("I4").Value2 = ''
$ws.Range("J4").Value2 = 'Impact Journals'
$ws.Range("E5").Value2 = '[Robert%De Smet%NULL%0, Bea%Mellaerts%NULL%1, Hannelore%Vandewinckele%NULL%1, Peter%Lybeert%NULL%1, Eric%Frans%NULL%1, Sara%Ombelet%NULL%1, Wim%Lemahieu%NULL%1, Rolf%Symons%NULL%0, Erwin%Ho%NULL%1, Johan%Frans%NULL%0, Annick%Smismans%NULL%1, Michaël R.%Laurent%NULL%1]'
$ws.Range("I5").Value2 = ''
$ws.Range("J5").Value2 = 'AMDA - The Society for Post-Acute and Long-Term Care Medicine.'
$ws.Range("E6").Value2 = '[Marcello%Covino%marcello.covino@policlinicogemelli.it%0, Giuseppe%De Matteis%NULL%3, Giuseppe%De Matteis%NULL%0, Michele%Santoro%NULL%2, Luca%Sabia%NULL%2, Benedetta%Simeoni%NULL%1, Marcello%Candelli%NULL%1, Veronica%Ojetti%NULL%1, Francesco%Franceschi%NULL%2]'
$ws.Range("I6").Value2 = ''
$ws.Range("J6").Value2 = 'John Wiley &amp; Sons Australia, Ltd'
$ws.Range("D7").Value2 = 'Coronavirus disease 2019 (COVID‐19) is a global pandemic and information on risk factors for worse prognosis is needed to accurately identify patients at risk and potentially provide insight into therapeutic options.
 In this retrospective cohort study, including 3703 patients with laboratory confirmed COVID‐19, we identified risk factors associated with all‐cause mortality, need for hospitalization and mechanical ventilation.
 Male gender was independently associated with increased risk of hospitalization (adjusted odds ratio [ORadj]: 1.62; 95% confidence interval [95% CI]: 1.38‐1.91)), mechanical ventilation (ORadj: 1.35; 95% CI: 1.08‐1.69) and death (ORadj: 1.46; 95% CI: 1.17‐1.82).
 Patients &gt; 60 years had higher risk of hospitalization (ORadj: 5.47; 95% CI: 4.29‐6.96), mechanical ventilation (ORadj: 3.26; 95% CI: 2.08‐5.11) and death (ORadj: 13.04; 95% CI: 6.25‐27.24).
 Congestive heart failure (ORadj: 1.47; 95% CI: 1.06‐2.02) and dementia (ORadj: 2.03; 95% CI: 1.46‐2.83) were associated with increased odds of death, as well as the presence of more than two comorbidities (ORadj: 1.90; 95% CI: 1.35‐2.68).
 Patients with COVID‐19 of older age, male gender, or having more than two comorbidities are at higher risk of hospitalization, mechanical ventilation and death, and should therefore be closely monitored.
'
$ws.Range("E7").Value2 = '[Maaike%van Gerwen%maaike.vangerwen@mountsinai.org%0, Mathilda%Alsen%NULL%2, Mathilda%Alsen%NULL%0, Christine%Little%NULL%1, Joshua%Barlow%NULL%1, Eric%Genden%NULL%1, Leonard%Naymagon%NULL%1, Douglas%Tremblay%NULL%1]'
$ws.Range("I7").Value2 = ''
$ws.Range("J7").Value2 = 'John Wiley and Sons Inc.'
$ws.Range("D8").Value2 = 'Introduction
id="Par1">In the current study, we evaluated factors that increase the coronavirus disease (COVID-19) patient death rate by analyzing the data from two cohort hospitals.

 In addition, we studied whether underlying neurological diseases are risk factors for death.


Methods
id="Par2">In this retrospective cohort study, we included 103 adult inpatients (aged ≥ 18 years).

 We evaluated differences in demographic data between surviving and non-surviving COVID-19 patients.


Results
id="Par3">In a multivariate logistic analysis, age and the presence of chronic lung disease and Alzheimer’s dementia (AD) were the only significant parameters for predicting COVID-19 non-survival (p &lt; 0.05).

 However, hypertension, coronary vascular disease, dyslipidemia, chronic kidney disease, diabetes, and history of taking angiotensin II receptor blockers (ARBs) or angiotensin-converting enzyme (ACE) inhibitors, as well as nonsteroidal anti-inflammatory drugs (NSAIDs), were not significantly associated with the death of COVID-19 patients.

 The optimal cutoff value obtained from the maximum Youden index was 70 (sensitivity, 80.77%; specificity, 61.04%), and the odds ratio of non-survival increased 1.055 fold for every year of age.


Conclusions
id="Par4">Clinicians should closely monitor and manage the symptoms of COVID-19 patients who are over the age of 70 years or have chronic lung disease or AD.


'
$ws.Range("E8").Value2 = '[Jong-moon%Hwang%NULL%0, Ju-Hyun%Kim%NULL%0, Jin-Sung%Park%NULL%0, Min Cheol%Chang%wheel633@ynu.ac.kr%0, Donghwi%Park%bdome@hanmail.net%0]'
$ws.Range("I8").Value2 = ''
$ws.Range("J8").Value2 = 'Springer International Publishing'
$ws.Range("D9").Value2 = 'id="Par1">The natural history of coronavirus disease 2019 (COVID-19) has yet to be fully described.
 Here, we use patient-level data from the Information System for Research in Primary Care (SIDIAP) to summarise COVID-19 outcomes in Catalonia, Spain.
 We included 5,586,521 individuals from the general population.
 Of these, 102,002 had an outpatient diagnosis of COVID-19, 16,901 were hospitalised with COVID-19, and 5273 died after either being diagnosed or hospitalised with COVID-19 between 1st March and 6th May 2020. Older age, being male, and having comorbidities were all generally associated with worse outcomes.
 These findings demonstrate the continued need to protect those at high risk of poor outcomes, particularly older people, from COVID-19 and provide appropriate care for those who develop symptomatic disease.
 While risks of hospitalisation and death were lower for younger populations, there is a need to limit their role in community transmission.
'
$ws.Range("E9").Value2 = '[Edward%Burn%NULL%0, Cristian%Tebé%NULL%2, Cristian%Tebé%NULL%0, Sergio%Fernandez-Bertolin%NULL%1, Maria%Aragon%NULL%1, Martina%Recalde%NULL%2, Martina%Recalde%NULL%0, Elena%Roel%NULL%1, Albert%Prats-Uribe%NULL%3, Albert%Prats-Uribe%NULL%0, Daniel%Prieto-Alhambra%daniel.prietoalhambra@ndorms.ox.ac.uk%3, Daniel%Prieto-Alhambra%daniel.prietoalhambra@ndorms.ox.ac.uk%0, Talita%Duarte-Salles%NULL%2, Talita%Duarte-Salles%NULL%0]'
$ws.Range("I9").Value2 = ''
$ws.Range("J9").Value2 = 'Nature Publishing Group UK'
$ws.Range("D10").Value2 = 'This is a population-based prospective cohort study on archive data describing the age- and sex-specific prevalence of COVID-19 and its prognostic factors.
 All 2653 symptomatic patients tested positive for SARS-CoV-2 from February 27 to April 2, 2020 in the Reggio Emilia province, Italy, were included.
 COVID-19 cumulative incidence, hospitalization and death rates, and adjusted hazard ratios (HR) with 95% confidence interval (95% CI) were calculated according to sociodemographic and clinical characteristics.
 Females had higher prevalence of infection than males below age 50 (2.61 vs.
 1.84 ‰), but lower in older ages (16.49 vs.
 20.86 ‰ over age 80).
 Case fatality rate reached 20.7% in cases with more than 4 weeks follow up.
 After adjusting for age and comorbidities, men had a higher risk of hospitalization (HR 1.4 95% CI 1.2 to 1.6) and of death (HR 1.6, 95% CI 1.2 to 2.1).
 Patients over age 80 compared to age &lt; 50 had HR 7.1 (95% CI 5.4 to 9.3) and HR 27.8 (95% CI 12.5 to 61.7) for hospitalization and death, respectively.
 Immigrants had a higher risk of hospitalization (HR 1.3, 95% CI 0.99 to 1.81) than Italians and a similar risk of death.
 Risk of hospitalization and of death were higher in patients with heart failure, arrhythmia, dementia, coronary heart disease, diabetes, and hypertension, while COPD increased the risk of hospitalization (HR 1.9, 95% CI 1.4 to 2.5) but not of death (HR 1.1, 95% CI 0.7 to 1.7).
 Previous use of ACE inhibitors had no effect on risk of death (HR 0.97, 95% CI 0.69 to 1.34).
 Identified susceptible populations and fragile patients should be considered when setting priorities in public health planning and clinical decision making.
'
$ws.Range("E10").Value2 = '[Paolo%Giorgi Rossi%NULL%0, Massimiliano%Marino%NULL%0, Debora%Formisano%NULL%0, Francesco%Venturelli%NULL%0, Massimo%Vicentini%NULL%0, Massimo%Vicentini%NULL%0, Roberto%Grilli%NULL%0, NULL%NULL%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0, Gianluigi%Forloni%NULL%0]'
$ws.Range("I10").Value2 = ''
$ws.Range("J10").Value2 = 'Public Library of Science'
$ws.Range("E11").Value2 = '[Sanjana P.%Padala%NULL%0, Ashlyn M.%Jendro%NULL%1, Lillian C.%Orr%Lillian.Orr@va.gov%1]'
$ws.Range("I11").Value2 = ''
$ws.Range("J11").Value2 = 'Elsevier/North-Holland Biomedical Press'
